$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.480199999999999
$ws.Range("A8").Value = -22.30400000000002
$ws.Range("A10").Value = -22.0207
$ws.Range("D11").Value = -7.157899999999993
$ws.Range("A12").Value = -21.4515
$ws.Range("D12").Value = -5.854499999999996
$ws.Range("D15").Value = -8.471399999999996
$ws.Range("D17").Value = -8.222899999999992
$ws.Range("A18").Value = -22.032
$ws.Range("A25").Value = -21.8149
$ws.Range("D26").Value = -7.066900000000008
$ws.Range("D27").Value = -8.439700000000002
$ws.Range("D28").Value = -8.735099999999999
$ws.Range("D32").Value = -6.205999999999992
$ws.Range("A37").Value = -19.98859999999998
$ws.Range("D37").Value = -7.903100000000002
$ws.Range("D41").Value = -8.272099999999993
$ws.Range("D47").Value = -7.900599999999998
$ws.Range("D51").Value = -8.237399999999999
$ws.Range("A55").Value = -21.888
$ws.Range("D65").Value = -7.842300000000002
$ws.Range("A68").Value = -21.39109999999999
$ws.Range("D73").Value = -8.414699999999995
$ws.Range("A77").Value = -19.93249999999999
$ws.Range("A78").Value = -19.56209999999998
$ws.Range("A79").Value = -20.16819999999998
$ws.Range("A80").Value = -19.44290000000001
$ws.Range("A81").Value = -22.0885
$ws.Range("A82").Value = -21.5415
$ws.Range("A84").Value = -21.9913
$ws.Range("D84").Value = -8.269500000000003
$ws.Range("D85").Value = -8.183499999999999
$ws.Range("D89").Value = -8.344899999999997
$ws.Range("D93").Value = -6.479899999999991
$ws.Range("D95").Value = -7.579200000000003
$ws.Range("D98").Value = -6.976400000000003
$ws.Range("D99").Value = -8.194400000000007
$ws.Range("A101").Value = -20.68799999999999
$ws.Range("D101").Value = -8.253699999999995
$ws.Range("A102").Value = -19.72579999999998
$ws.Range("D102").Value = -7.8728